# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns for both locales now that the handback round-trip completed, and refreshes the
# "Status" text everywhere it is echoed (Overview + each locale sheet) to reflect the
# handed-back / in-sync state.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$zhCnLink1 = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6af5b141368fd80643c2c0a33704092effd2456/e2e/362ff2f9-3b05-40e3-8c11-97be5941e3cb.md", "", "", "362ff2f9-3b05-40e3-8c11-97be5941e3cb.md")
$wsZhCn.Range("J2").Value = "362ff2f9-3b05-40e3-8c11-97be5941e3cb.9b53ffc75c5239980ee036c680d0002980ae0673.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-05 19:09:04"

$zhCnLink2 = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6af5b141368fd80643c2c0a33704092effd2456/e2e/ac767584-6af0-470d-a3c7-014cc8455e3f.md", "", "", "ac767584-6af0-470d-a3c7-014cc8455e3f.md")
$wsZhCn.Range("J3").Value = "ac767584-6af0-470d-a3c7-014cc8455e3f.f4ee9877fabf1f4eec22d587239436f278c4fc2e.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-05 19:09:04"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

$deDeLink1 = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6af5b141368fd80643c2c0a33704092effd2456/e2e/362ff2f9-3b05-40e3-8c11-97be5941e3cb.md", "", "", "362ff2f9-3b05-40e3-8c11-97be5941e3cb.md")
$wsDeDe.Range("J2").Value = "362ff2f9-3b05-40e3-8c11-97be5941e3cb.9b53ffc75c5239980ee036c680d0002980ae0673.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-05 19:09:23"

$deDeLink2 = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6af5b141368fd80643c2c0a33704092effd2456/e2e/ac767584-6af0-470d-a3c7-014cc8455e3f.md", "", "", "ac767584-6af0-470d-a3c7-014cc8455e3f.md")
$wsDeDe.Range("J3").Value = "ac767584-6af0-470d-a3c7-014cc8455e3f.f4ee9877fabf1f4eec22d587239436f278c4fc2e.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-05 19:09:23"

# Widen the columns that now hold longer text so everything stays readable,
# mirroring the autofit Excel performs after long values are dropped in.
$wsOverview.Range("E1:F1").ColumnWidth = 29.9777047293527

$wsZhCn.Range("C1").ColumnWidth = 29.9777047293527
$wsZhCn.Range("I1:J1").ColumnWidth = 40

$wsDeDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDeDe.Range("I1:J1").ColumnWidth = 40
